$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 5

$ws.Range("E2").Value = "dsfdsfdsaf"
$ws.Range("E5").Value = "dsf"
$ws.Range("E6").Value = "dsfds"

$ws.Range("J17").Select()
